$wb = $excel.ActiveWorkbook

# --- Creational sheet updates ---
$wsCreational = $wb.Worksheets.Item("Creational")
$wsCreational.Range("D5").Value = "Chưa rõ"
$wsCreational.Range("D6").Value = "Chưa rõ"
$wsCreational.Range("C7").Value = "Việc khởi tạo thực hiện duy nhất 1 lần"
$wsCreational.Range("D7").Value = "Có thể áp dụng"

# --- Behavioral sheet updates ---
$wsBehavioral = $wb.Worksheets.Item("Behavioral")
$wsBehavioral.Range("C3").Value = "Quyết định quy trình chạy của các lớp theo thứ tự nhất định"
$wsBehavioral.Range("D3").Value = "Chưa rõ"
$wsBehavioral.Range("C4").Value = "Chuyên xử lý các hành động undo, redo"
$wsBehavioral.Range("D4").Value = "Chưa rõ"
$wsBehavioral.Range("C5").Value = "Quản lý định dạng date hoặc việc đọc dữ liệu từ bên ngoài của nhiều class"
$wsBehavioral.Range("D5").Value = "Chưa rõ"
$wsBehavioral.Range("C6").Value = "Quản lý kiểu danh sách nhiều phần tử"
$wsBehavioral.Range("D6").Value = "Chưa rõ"
$wsBehavioral.Range("C7").Value = "Điều phối thông điệp với các thể hiện khác nhau"
$wsBehavioral.Range("D7").Value = "Chưa rõ"

# Make Behavioral the active sheet, with C7 selected as the active cell.
$wsBehavioral.Activate()
$wsBehavioral.Range("C7").Select()
